# Daily backup 2025-11-30 01:26:07
# Update column D values in the data table and refresh the saved view/selection state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("D319").Value = 3
$ws.Range("D330").Value = 3
$ws.Range("D341").Value = 3
$ws.Range("D345").Value = 43.531999999999996
$ws.Range("D347").Value = 44
$ws.Range("D352").Value = 3
$ws.Range("D356").Value = 43.012
$ws.Range("D357").Value = 11.792
$ws.Range("D358").Value = 42
$ws.Range("D363").Value = 3
$ws.Range("D367").Value = 42.472000000000001
$ws.Range("D368").Value = 12.2005
$ws.Range("D369").Value = 48
$ws.Range("D374").Value = 3
$ws.Range("D378").Value = 42.091999999999999
$ws.Range("D379").Value = 12.659000000000001
$ws.Range("D380").Value = 46
$ws.Range("D385").Value = 3
$ws.Range("D389").Value = 39.771999999999998
$ws.Range("D390").Value = 10.366
$ws.Range("D391").Value = 43
$ws.Range("D396").Value = 3
$ws.Range("D400").Value = 39.392000000000003
$ws.Range("D401").Value = 10.683999999999999
$ws.Range("D402").Value = 34
$ws.Range("D407").Value = 3
$ws.Range("D411").Value = 39.012
$ws.Range("D412").Value = 11.007
$ws.Range("D413").Value = 38
$ws.Range("D418").Value = 3
$ws.Range("D422").Value = 38.631999999999998
$ws.Range("D423").Value = 11.3605
$ws.Range("D424").Value = 44
$ws.Range("D429").Value = 3
$ws.Range("D433").Value = 38.351999999999997
$ws.Range("D434").Value = 11.743
$ws.Range("D435").Value = 36
$ws.Range("D440").Value = 3
$ws.Range("D444").Value = 38.131999999999998
$ws.Range("D445").Value = 12.157
$ws.Range("D446").Value = 43
$ws.Range("D451").Value = 3
$ws.Range("D455").Value = 37.771999999999998
$ws.Range("D456").Value = 12.608000000000001
$ws.Range("D457").Value = 39
$ws.Range("D462").Value = 3
$ws.Range("D473").Value = 3
$ws.Range("D484").Value = 3
$ws.Range("D495").Value = 3
$ws.Range("D506").Value = 3
$ws.Range("D517").Value = 3
$ws.Range("D528").Value = 3
$ws.Range("D539").Value = 3

# Restore the view state (scroll position + active selection) recorded in the workbook.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 463
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G471").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("B463")
